$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 721.7
$ws.Cells.Item(12, 9).Value = 513.4286
$ws.Cells.Item(12, 10).Value = 1207.6666
$ws.Cells.Item(12, 11).Value = 513.4286
$ws.Cells.Item(12, 12).Value = 1207.6666
$ws.Cells.Item(12, 13).Value = -343.4286
$ws.Cells.Item(12, 14).Value = -1547.6666
$ws.Cells.Item(51, 8).Value = 5635.909
$ws.Cells.Item(51, 10).Value = 5600
$ws.Cells.Item(51, 12).Value = 5600
$ws.Cells.Item(51, 14).Value = -6568
$ws.Cells.Item(111, 8).Value = 1382.4286
$ws.Cells.Item(111, 9).Value = 1642.6666
$ws.Cells.Item(111, 10).Value = 1187.25
$ws.Cells.Item(111, 11).Value = 4927.9998
$ws.Cells.Item(111, 12).Value = 3561.75
$ws.Cells.Item(111, 13).Value = -1860.9998
$ws.Cells.Item(111, 14).Value = -9695.75
$ws.Cells.Item(116, 8).Value = 7722.909
$ws.Cells.Item(116, 9).Value = 7051.8335
$ws.Cells.Item(116, 11).Value = 7051.8335
$ws.Cells.Item(116, 13).Value = -3609.8335
$ws.Cells.Item(127, 8).Value = 1183
$ws.Cells.Item(127, 9).Value = 1118.9474
$ws.Cells.Item(127, 11).Value = 3356.8422
$ws.Cells.Item(127, 13).Value = 1603.1578
$ws.Cells.Item(141, 8).Value = 4886.6
$ws.Cells.Item(141, 9).Value = 5221.625
$ws.Cells.Item(141, 11).Value = 15664.875
$ws.Cells.Item(141, 13).Value = -10484.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 37037640
$ws.Cells.Item(74, 9).Value = 37037640
$ws.Cells.Item(74, 11).Value = 37037640
$ws.Cells.Item(74, 13).Value = -37036766
$ws.Cells.Item(77, 8).Value = 37037640
$ws.Cells.Item(77, 9).Value = 37037640
$ws.Cells.Item(77, 11).Value = 185188200
$ws.Cells.Item(77, 13).Value = -185183832
$ws.Cells.Item(97, 8).Value = 3261.1428
$ws.Cells.Item(97, 10).Value = 5771.4287
$ws.Cells.Item(97, 12).Value = 5771.4287
$ws.Cells.Item(97, 14).Value = -6763.4287
$ws.Cells.Item(133, 8).Value = 63980.5
$ws.Cells.Item(133, 10).Value = 63980.5
$ws.Cells.Item(133, 12).Value = 63980.5
$ws.Cells.Item(133, 14).Value = -69040.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 342.13043
$ws.Cells.Item(22, 9).Value = 328.57144
$ws.Cells.Item(22, 11).Value = 328.57144
$ws.Cells.Item(22, 13).Value = -155.57144
$ws.Cells.Item(63, 8).Value = 63333.332
$ws.Cells.Item(63, 10).Value = 63333.332
$ws.Cells.Item(63, 12).Value = 63333.332
$ws.Cells.Item(63, 14).Value = -64705.332
$ws.Cells.Item(66, 8).Value = 63333.332
$ws.Cells.Item(66, 10).Value = 63333.332
$ws.Cells.Item(66, 12).Value = 189999.996
$ws.Cells.Item(66, 14).Value = -196863.996
$ws.Cells.Item(94, 8).Value = 1266.3793
$ws.Cells.Item(94, 9).Value = 1294.5454
$ws.Cells.Item(94, 11).Value = 1294.5454
$ws.Cells.Item(94, 13).Value = -843.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1912.6
$ws.Cells.Item(16, 10).Value = 1674.5
$ws.Cells.Item(16, 12).Value = 1674.5
$ws.Cells.Item(16, 14).Value = -2248.5
$ws.Cells.Item(22, 8).Value = 925.6070999999999
$ws.Cells.Item(22, 9).Value = 693
$ws.Cells.Item(22, 10).Value = 1035.7894
$ws.Cells.Item(22, 11).Value = 693
$ws.Cells.Item(22, 12).Value = 1035.7894
$ws.Cells.Item(22, 13).Value = -343
$ws.Cells.Item(22, 14).Value = -1735.7894
$ws.Cells.Item(58, 8).Value = 1395.2142
$ws.Cells.Item(58, 9).Value = 1412.7142
$ws.Cells.Item(58, 10).Value = 1377.7142
$ws.Cells.Item(58, 11).Value = 1412.7142
$ws.Cells.Item(58, 12).Value = 1377.7142
$ws.Cells.Item(58, 13).Value = -1209.7142
$ws.Cells.Item(58, 14).Value = -1783.7142
$ws.Cells.Item(113, 8).Value = 1912.6
$ws.Cells.Item(113, 10).Value = 1674.5
$ws.Cells.Item(113, 12).Value = 1674.5
$ws.Cells.Item(113, 14).Value = -6014.5
$ws.Cells.Item(132, 8).Value = 2408
$ws.Cells.Item(132, 9).Value = 2311.4062
$ws.Cells.Item(132, 11).Value = 6934.2186
$ws.Cells.Item(132, 13).Value = -4404.2186
$ws.Cells.Item(134, 8).Value = 4259.4346
$ws.Cells.Item(134, 9).Value = 4291
$ws.Cells.Item(134, 10).Value = 4145.8
$ws.Cells.Item(134, 11).Value = 12873
$ws.Cells.Item(134, 12).Value = 12437.4
$ws.Cells.Item(134, 13).Value = -10338
$ws.Cells.Item(134, 14).Value = -17507.4
$ws.Cells.Item(135, 8).Value = 70567.38
$ws.Cells.Item(135, 10).Value = 70567.38
$ws.Cells.Item(135, 12).Value = 70567.38
$ws.Cells.Item(135, 14).Value = -80707.38
$ws.Cells.Item(136, 8).Value = 1395.2142
$ws.Cells.Item(136, 9).Value = 1412.7142
$ws.Cells.Item(136, 10).Value = 1377.7142
$ws.Cells.Item(136, 11).Value = 4238.142599999999
$ws.Cells.Item(136, 12).Value = 4133.142599999999
$ws.Cells.Item(136, 13).Value = -1688.142599999999
$ws.Cells.Item(136, 14).Value = -9233.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 319.44827
$ws.Cells.Item(12, 10).Value = 343.23077
$ws.Cells.Item(12, 12).Value = 1029.69231
$ws.Cells.Item(12, 14).Value = -1375.69231
$ws.Cells.Item(14, 8).Value = 2099.1428
$ws.Cells.Item(14, 9).Value = 2099.1428
$ws.Cells.Item(14, 11).Value = 6297.428400000001
$ws.Cells.Item(14, 13).Value = -6124.428400000001
$ws.Cells.Item(68, 8).Value = 1384.5714
$ws.Cells.Item(68, 9).Value = 1088.8667
$ws.Cells.Item(68, 10).Value = 1606.35
$ws.Cells.Item(68, 11).Value = 3266.6001
$ws.Cells.Item(68, 12).Value = 4819.049999999999
$ws.Cells.Item(68, 13).Value = -2455.6001
$ws.Cells.Item(68, 14).Value = -6441.049999999999
$ws.Cells.Item(71, 8).Value = 1384.5714
$ws.Cells.Item(71, 9).Value = 1088.8667
$ws.Cells.Item(71, 10).Value = 1606.35
$ws.Cells.Item(71, 11).Value = 9799.800300000001
$ws.Cells.Item(71, 12).Value = 14457.15
$ws.Cells.Item(71, 13).Value = -5743.800300000001
$ws.Cells.Item(71, 14).Value = -22569.15
$ws.Cells.Item(80, 8).Value = 29999.6
$ws.Cells.Item(80, 10).Value = 20000
$ws.Cells.Item(80, 12).Value = 60000
$ws.Cells.Item(80, 14).Value = -61872
$ws.Cells.Item(83, 8).Value = 29999.6
$ws.Cells.Item(83, 10).Value = 20000
$ws.Cells.Item(83, 12).Value = 180000
$ws.Cells.Item(83, 14).Value = -189360
$ws.Cells.Item(136, 8).Value = 4472.5
$ws.Cells.Item(136, 9).Value = 4472.5
$ws.Cells.Item(136, 11).Value = 13417.5
$ws.Cells.Item(136, 13).Value = -8317.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 6689.421
$ws.Cells.Item(102, 9).Value = 6803.643
$ws.Cells.Item(102, 10).Value = 6369.6
$ws.Cells.Item(102, 11).Value = 6803.643
$ws.Cells.Item(102, 12).Value = 6369.6
$ws.Cells.Item(102, 13).Value = -5181.643
$ws.Cells.Item(102, 14).Value = -9613.6
$ws.Cells.Item(113, 8).Value = 37044344
$ws.Cells.Item(113, 10).Value = 4099.857
$ws.Cells.Item(113, 12).Value = 4099.857
$ws.Cells.Item(113, 14).Value = -8439.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2809.5557
$ws.Cells.Item(22, 9).Value = 2047.6666
$ws.Cells.Item(22, 10).Value = 4333.3335
$ws.Cells.Item(22, 11).Value = 2047.6666
$ws.Cells.Item(22, 12).Value = 4333.3335
$ws.Cells.Item(22, 13).Value = -1752.6666
$ws.Cells.Item(22, 14).Value = -4923.3335
$ws.Cells.Item(27, 8).Value = 2809.5557
$ws.Cells.Item(27, 9).Value = 2047.6666
$ws.Cells.Item(27, 10).Value = 4333.3335
$ws.Cells.Item(27, 11).Value = 2047.6666
$ws.Cells.Item(27, 12).Value = 4333.3335
$ws.Cells.Item(27, 13).Value = -1940.6666
$ws.Cells.Item(27, 14).Value = -4547.3335
$ws.Cells.Item(40, 8).Value = 6146.2856
$ws.Cells.Item(40, 9).Value = 5809.091
$ws.Cells.Item(40, 11).Value = 5809.091
$ws.Cells.Item(40, 13).Value = -5673.091
$ws.Cells.Item(132, 8).Value = 2397152.5
$ws.Cells.Item(132, 9).Value = 2529494.5
$ws.Cells.Item(132, 11).Value = 7588483.5
$ws.Cells.Item(132, 13).Value = -7585953.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 10205629
$ws.Cells.Item(81, 9).Value = 17858940
$ws.Cells.Item(81, 10).Value = 1215
$ws.Cells.Item(81, 11).Value = 35717880
$ws.Cells.Item(81, 12).Value = 2430
$ws.Cells.Item(81, 13).Value = -35716819
$ws.Cells.Item(81, 14).Value = -4552
$ws.Cells.Item(84, 8).Value = 10205629
$ws.Cells.Item(84, 9).Value = 17858940
$ws.Cells.Item(84, 10).Value = 1215
$ws.Cells.Item(84, 11).Value = 178589400
$ws.Cells.Item(84, 12).Value = 12150
$ws.Cells.Item(84, 13).Value = -178584096
$ws.Cells.Item(84, 14).Value = -22758
$ws.Cells.Item(113, 8).Value = 2532.5
$ws.Cells.Item(113, 9).Value = 1359.2858
$ws.Cells.Item(113, 11).Value = 4077.8574
$ws.Cells.Item(113, 13).Value = -1907.8574
$ws.Cells.Item(119, 8).Value = 92349.64999999999
$ws.Cells.Item(119, 10).Value = 92349.64999999999
$ws.Cells.Item(119, 12).Value = 92349.64999999999
$ws.Cells.Item(119, 14).Value = -102025.65
$ws.Cells.Item(126, 8).Value = 113758450
$ws.Cells.Item(126, 9).Value = 113758450
$ws.Cells.Item(126, 11).Value = 341275350
$ws.Cells.Item(126, 13).Value = -341272880
$ws.Cells.Item(132, 8).Value = 3414.4119
$ws.Cells.Item(132, 9).Value = 3424.6875
$ws.Cells.Item(132, 11).Value = 10274.0625
$ws.Cells.Item(132, 13).Value = -7744.0625
$ws.Cells.Item(136, 8).Value = 3498292.8
$ws.Cells.Item(136, 9).Value = 5496075
$ws.Cells.Item(136, 10).Value = 2174.375
$ws.Cells.Item(136, 11).Value = 16488225
$ws.Cells.Item(136, 12).Value = 6523.125
$ws.Cells.Item(136, 13).Value = -16485675
$ws.Cells.Item(136, 14).Value = -11623.125
